$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new ingredient row (amora / blackberry pulp)
$ws.Range("A37").Value = "amora"
$ws.Range("B37").Value = "Amora (polpa)"
$ws.Range("C37").Value = "fruit"
$ws.Range("D37").Value = 0.2
$ws.Range("E37").Value = 11.7
$ws.Range("F37").Value = 13
$ws.Range("G37").Value = 1.02
$ws.Range("H37").Value = $true
$ws.Range("J37").Value = $true
$ws.Range("K37").Value = "Brix ~13; enzimas podem afetar textura."
$ws.Range("L37").Value = "1.0.0"
$ws.Range("O37").Value = "OK"

# Update view state to match the author's saved selection
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("G37:O37").Select()
